$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.437.25"
$ws.Range("E2").Value = "  -5.91%  "

$ws.Range("D3").Value = "2.926.26"
$ws.Range("E3").Value = "  -9.40%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "469.74"
$ws.Range("E5").Value = "  -12.57%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "123.16"
$ws.Range("E6").Value = "  -9.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "2.927.19"
$ws.Range("E8").Value = "  -9.27%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.392"
$ws.Range("E9").Value = "  -14.46%  "

$ws.Range("E10").Value = "  -13.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0943"
$ws.Range("E11").Value = "  -18.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.326"
$ws.Range("E12").Value = "  -17.56%  "

$ws.Range("E13").Value = "  -3.31%  "

$ws.Range("D14").Value = "3.434.02"
$ws.Range("E14").Value = "  -9.10%  "

$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.38"
$ws.Range("E15").Value = "  -14.10%  "

$ws.Range("D16").Value = "55.630.49"
$ws.Range("E16").Value = "  -5.71%  "

$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.960.26"
$ws.Range("E17").Value = "  -8.31%  "

$ws.Range("E18").Value = "  -17.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.07"
$ws.Range("E19").Value = "  -14.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.41"
$ws.Range("E20").Value = "  -13.84%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.99"
$ws.Range("E21").Value = "  -15.59%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "306.10"
$ws.Range("E22").Value = "  -15.28%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.443"
$ws.Range("E24").Value = "  -14.71%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "58.72"
$ws.Range("E25").Value = "  -16.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.154"
$ws.Range("E27").Value = "  -9.71%  "

$ws.Range("E28").Value = "  -0.13%  "

$ws.Range("D29").Value = "0.0₃0794"
$ws.Range("E29").Value = "  -18.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.86"
$ws.Range("E30").Value = "  -17.25%  "

$ws.Range("E31").Value = "  -10.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.99"
$ws.Range("E32").Value = "  -13.80%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.00"
$ws.Range("E33").Value = "  -15.45%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.56"
$ws.Range("E34").Value = "  -19.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "145.07"
$ws.Range("E35").Value = "  -11.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.14"
$ws.Range("E36").Value = "  -16.19%  "

$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.20"
$ws.Range("E37").Value = "  -16.38%  "

$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.30"
$ws.Range("E38").Value = "  -16.81%  "

$ws.Range("B39").Value = "FirstDigitalUSD"
$ws.Range("C39").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.30%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0603"
$ws.Range("E40").Value = "  -14.69%  "

$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "21.24"
$ws.Range("E41").Value = "  -18.21%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "34.68"
$ws.Range("E42").Value = "  -15.50%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.597"
$ws.Range("E43").Value = "  -16.69%  "

$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.936"
$ws.Range("E44").Value = "  -14.03%  "

$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.37"
$ws.Range("E45").Value = "  -16.07%  "

$ws.Range("D46").Value = "2.056.27"
$ws.Range("E46").Value = "  -10.25%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.29"
$ws.Range("E47").Value = "  -14.33%  "

$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.21"
$ws.Range("E48").Value = "  -16.78%  "

$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.38"
$ws.Range("E49").Value = "  -16.13%  "

$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0208"
$ws.Range("E50").Value = "  -14.16%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0786"
$ws.Range("E51").Value = "  -11.33%  "
